$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of JPYUSD FX data to append starting at row 882
# Columns: A=datetime(serial), B=open, C=high, D=low, E=close, F=volume
$data = @(
    @(45260, 0.00679,   0.006809,  0.006736,  0.006745, 0),
    @(45261, 0.006746,  0.006817,  0.006743,  0.006809, 0),
    @(45261, 0.006746,  0.006817,  0.006743,  0.006809, 0),
    @(45261, 0.006746,  0.006817,  0.006743,  0.006809, 0),
    @(45264, 0.00681,   0.006838,  0.006783,  0.006791, 0),
    @(45265, 0.006792,  0.006821,  0.006785,  0.006794, 0),
    @(45266, 0.006794,  0.0067979, 0.0067856, 0.006794, 0),
    @(45266, 0.006794,  0.006806,  0.006781,  0.006786, 0),
    @(45267, 0.00679,   0.007049,  0.006786,  0.006934, 0),
    @(45268, 0.006933,  0.0070154, 0.0069301, 0.006937, 0)
)

$startRow = 881 + 1
$rowCount = $data.Count
$endRow = $startRow + $rowCount - 1

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# Reuse the existing date/time cell style from the last pre-existing row (A881)
# rather than building up a brand new style, matching how the rest of column A
# is formatted (YYYY-MM-DD HH:MM:SS, centered, bordered).
$ws.Range("A881").Copy()
$ws.Range("A$startRow`:A$endRow").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false | Out-Null
